# ng_oncho_2410_4_b_flies_lab_tar -> _v2
# Rename pools_positive/pools_negative/pools_less_100 fields (and all their
# dependents/references) to pools_pos/pools_neg/pools_l_100, bump the form
# id/title to the V2 variant, and move the active sheet/selection the way
# the authored workbook ended up.

$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# --- "positive" pools repeat group: pools_positive* -> pools_pos* ---
$survey.Range("B10").Value = "pools_pos"
$survey.Range("B11").Value = "pools_pos_details_tar"
$survey.Range("H11").Value = '${pools_pos} > 0'
$survey.Range("N11").Value = '${pools_pos}'
$survey.Range("B12").Value = "pools_pos_name"
$survey.Range("B13").Value = "pools_pos_result"
$survey.Range("B14").Value = "pools_pos_size"
$survey.Range("B15").Value = "pools_pos_test_type"
$survey.Range("C15").Value = 'Test type for pool **${pools_pos_name}**'

# --- "negative" pools repeat group: pools_negative* -> pools_neg* ---
$survey.Range("B17").Value = "pools_neg"
$survey.Range("B18").Value = "pools_neg_details_tar"
$survey.Range("H18").Value = '${pools_neg} > 0'
$survey.Range("N18").Value = '${pools_neg}'
$survey.Range("B19").Value = "pools_neg_name"
$survey.Range("B20").Value = "pools_neg_result"
$survey.Range("B21").Value = "pools_neg_size"
$survey.Range("B22").Value = "pools_neg_test_type"
$survey.Range("C22").Value = 'Test type for pool **${pools_neg_name}**'

# --- "less than 100" pools repeat group: pools_less_100* -> pools_l_100* ---
$survey.Range("B24").Value = "pools_l_100"
$survey.Range("B25").Value = "pools_l_100_details_tar"
$survey.Range("H25").Value = '${pools_l_100} > 0'
$survey.Range("N25").Value = '${pools_l_100}'
$survey.Range("B26").Value = "pools_l_100_name"
$survey.Range("B27").Value = "pools_l_100_test_type"
$survey.Range("C27").Value = 'Test type for pool **${pools_l_100_name}**'

# --- calculate rows referencing the renamed fields ---
$survey.Range("I33").Value = 'sum(${pool_result_negative}) + ${pools_neg}'
$survey.Range("I34").Value = 'sum(${pool_result_positive}) + ${pools_pos}'

# --- settings sheet: bump form_id / form_title to the V2 variant ---
$settings.Range("B2").Value = "ng_oncho_2410_4_b_flies_lab_tar_v2"
$settings.Range("A2").Value = "(Taraba) 4. Blackfly Lab App V2"

# --- selection / active sheet bookkeeping, matching the saved workbook ---
$survey.Range("B24").Select()
$settings.Activate()
